# Scheduled-runner update to the Leve profit tables (currentAveragePrice.. /
# LevePrice.. / LeveProfit.. columns H:N) across several sheets, reflecting
# refreshed market-board pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2725.25
$ws.Range("I62").Value = 3035
$ws.Range("J62").Value = 1796
$ws.Range("K62").Value = 3035
$ws.Range("L62").Value = 1796
$ws.Range("M62").Value = -2411
$ws.Range("N62").Value = -3044

$ws.Range("H65").Value = 2725.25
$ws.Range("I65").Value = 3035
$ws.Range("J65").Value = 1796
$ws.Range("K65").Value = 15175
$ws.Range("L65").Value = 8980
$ws.Range("M65").Value = -12055
$ws.Range("N65").Value = -15220

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 47619156
$ws.Range("I5").Value = 20833432
$ws.Range("J5").Value = 83333450
$ws.Range("K5").Value = 20833432
$ws.Range("L5").Value = 83333450
$ws.Range("M5").Value = -20833320
$ws.Range("N5").Value = -83333674

$ws.Range("H32").Value = 18452.12
$ws.Range("I32").Value = 18095.06
$ws.Range("J32").Value = 21380
$ws.Range("K32").Value = 18095.06
$ws.Range("L32").Value = 21380
$ws.Range("M32").Value = -17808.06
$ws.Range("N32").Value = -21954

$ws.Range("H109").Value = 20500
$ws.Range("J109").Value = 20500
$ws.Range("L109").Value = 20500
$ws.Range("N109").Value = -23274

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 47619156
$ws.Range("I4").Value = 20833432
$ws.Range("J4").Value = 83333450
$ws.Range("K4").Value = 20833432
$ws.Range("L4").Value = 83333450
$ws.Range("M4").Value = -20833317
$ws.Range("N4").Value = -83333680

$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 218.72728
$ws.Range("I2").Value = 271.5625
$ws.Range("J2").Value = 77.833336
$ws.Range("K2").Value = 1629.375
$ws.Range("L2").Value = 467.000016
$ws.Range("M2").Value = -1516.375
$ws.Range("N2").Value = -693.000016

$ws.Range("H3").Value = 3463.1482
$ws.Range("I3").Value = 1906.5625
$ws.Range("J3").Value = 5727.273
$ws.Range("K3").Value = 5719.6875
$ws.Range("L3").Value = 17181.819
$ws.Range("M3").Value = -5607.6875
$ws.Range("N3").Value = -17405.819

$ws.Range("H4").Value = 177.77777
$ws.Range("I4").Value = 177.77777
$ws.Range("K4").Value = 533.33331
$ws.Range("M4").Value = -421.33331

$ws.Range("H6").Value = 241.2
$ws.Range("I6").Value = 51.5
$ws.Range("K6").Value = 154.5
$ws.Range("M6").Value = -41.5

$ws.Range("H7").Value = 73.75
$ws.Range("I7").Value = 70
$ws.Range("J7").Value = 75
$ws.Range("K7").Value = 210
$ws.Range("L7").Value = 225
$ws.Range("M7").Value = -98
$ws.Range("N7").Value = -449

$ws.Range("H9").Value = 23817028
$ws.Range("I9").Value = 3000
$ws.Range("J9").Value = 25648878
$ws.Range("K9").Value = 9000
$ws.Range("L9").Value = 76946634
$ws.Range("M9").Value = -8776
$ws.Range("N9").Value = -76947082

$ws.Range("H10").Value = 178.85715
$ws.Range("I10").Value = 178.85715
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 536.5714499999999
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -397.5714499999999
$ws.Range("N10").ClearContents()

$ws.Range("H11").Value = 11113345
$ws.Range("I11").Value = 50.5
$ws.Range("J11").Value = 14288571
$ws.Range("K11").Value = 151.5
$ws.Range("L11").Value = 42865713
$ws.Range("M11").Value = -11.5
$ws.Range("N11").Value = -42865993

$ws.Range("H12").Value = 126.48
$ws.Range("I12").Value = 103.375
$ws.Range("J12").Value = 137.35294
$ws.Range("K12").Value = 310.125
$ws.Range("L12").Value = 412.05882
$ws.Range("M12").Value = -137.125
$ws.Range("N12").Value = -758.05882

$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

$ws.Range("H15").Value = 757
$ws.Range("I15").Value = 163.33333
$ws.Range("J15").Value = 979.625
$ws.Range("K15").Value = 489.99999
$ws.Range("L15").Value = 2938.875
$ws.Range("M15").Value = -349.99999
$ws.Range("N15").Value = -3218.875

$ws.Range("H21").Value = 348.92856
$ws.Range("I21").Value = 99
$ws.Range("J21").Value = 487.77777
$ws.Range("K21").Value = 297
$ws.Range("L21").Value = 1463.33331
$ws.Range("M21").Value = -124
$ws.Range("N21").Value = -1809.33331

$ws.Range("H26").Value = 16667689
$ws.Range("I26").Value = 67.5
$ws.Range("J26").Value = 25001500
$ws.Range("K26").Value = 202.5
$ws.Range("L26").Value = 75004500
$ws.Range("M26").Value = 85.5
$ws.Range("N26").Value = -75005076

$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H33").Value = 29713.666
$ws.Range("I33").Value = 29713.666
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 178281.996
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -177998.996
$ws.Range("N33").ClearContents()

$ws.Range("H40").Value = 320.0909
$ws.Range("I40").Value = 66.94444
$ws.Range("J40").Value = 623.86664
$ws.Range("K40").Value = 267.77776
$ws.Range("L40").Value = 2495.46656
$ws.Range("M40").Value = -198.77776
$ws.Range("N40").Value = -2633.46656

$ws.Range("H44").Value = 681.4
$ws.Range("I44").Value = 566.6667
$ws.Range("J44").Value = 730.5714
$ws.Range("K44").Value = 1700.0001
$ws.Range("L44").Value = 2191.7142
$ws.Range("M44").Value = -1302.0001
$ws.Range("N44").Value = -2987.7142

$ws.Range("H46").Value = 285.7143
$ws.Range("J46").Value = 800
$ws.Range("L46").Value = 2400
$ws.Range("N46").Value = -2582

$ws.Range("H51").Value = 699.5
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 699.5
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 2098.5
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -3018.5

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws.Range("H58").Value = 300
$ws.Range("I58").Value = 300
$ws.Range("K58").Value = 900
$ws.Range("M58").Value = -772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws.Range("H122").Value = 4452.3335
$ws.Range("I122").Value = 5604.273
$ws.Range("J122").Value = 2642.1428
$ws.Range("K122").Value = 16812.819
$ws.Range("L122").Value = 7926.428400000001
$ws.Range("M122").Value = -14362.819
$ws.Range("N122").Value = -12826.4284

$ws.Range("H132").Value = 5491.8125
$ws.Range("I132").Value = 7039.2256
$ws.Range("J132").Value = 2670.0588
$ws.Range("K132").Value = 21117.6768
$ws.Range("L132").Value = 8010.176399999999
$ws.Range("M132").Value = -18587.6768
$ws.Range("N132").Value = -13070.1764

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1464.5862
$ws.Range("I132").Value = 1279.8334
$ws.Range("J132").Value = 1949.5625
$ws.Range("K132").Value = 3839.5002
$ws.Range("L132").Value = 5848.6875
$ws.Range("M132").Value = -1309.5002
$ws.Range("N132").Value = -10908.6875
